# prep_server success - fix product_name typos in digital_warranty_code column (E)
#  - "LODON LIBERTY TOILE" -> "LONDON LIBERTY TOILE" (rows 34-37)
#  - "SOLID SUIT SKINNY TIE · GREY" -> "SOLID SUIT SLIM TIE · GREY" (rows 42-43,
#    these are actually the "Slim" tie SKU per columns B/C, label just had wrong size name)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E34").Value = "LONDON LIBERTY TOILE · SIZE S · BLACK"
$ws.Range("E35").Value = "LONDON LIBERTY TOILE · SIZE S · BLACK"
$ws.Range("E36").Value = "LONDON LIBERTY TOILE · SIZE L · BLACK"
$ws.Range("E37").Value = "LONDON LIBERTY TOILE · SIZE L · BLACK"

$ws.Range("E42").Value = "SOLID SUIT SLIM TIE · GREY"
$ws.Range("E43").Value = "SOLID SUIT SLIM TIE · GREY"

# Reflect the author's last-saved cursor position
$ws.Range("E15").Select() | Out-Null
